$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 52
    6  = 186
    14 = 1629
    15 = 49
    22 = 1414
    23 = 3347
    25 = 55
    27 = 1085
    32 = 53
    33 = 276
    34 = 400
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
